$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Described" shifts from D to E),
# making room for the new "Index" column.
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "Index"

# Data rows: Index cycles 1,2,3 for each (BodyType, BodyBaseType) triplet.
for ($r = 2; $r -le 37; $r++) {
    $idx = (($r - 2) % 3) + 1
    $ws.Cells.Item($r, 4).Value = $idx
}
